$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date value (serial 45202 = 2023-10-03) for every
# data row (2 through 440). Bump it by one day to 45203 (2023-10-04) for all rows.
for ($r = 2; $r -le 440; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
